# ---------------------------------------------------------------------------
# Weekly refresh of the Femacal de La Calera - Ciboulette dataset.
# The source feed re-shuffles which record (Fecha/Volumen/Precio columns)
# lands on which spreadsheet row each week, and appends one brand-new record.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Fecha" (D) / "Volumen" (J) values for every existing data row (2-257).
$djData = @(
  @(2,44274,120),
  @(3,44321,130),
  @(4,44463,160),
  @(5,44459,160),
  @(6,44488,150),
  @(7,44497,160),
  @(8,44425,160),
  @(9,44215,130),
  @(10,44407,160),
  @(11,44551,160),
  @(12,44284,180),
  @(13,44267,160),
  @(14,44344,160),
  @(15,44428,160),
  @(16,44168,160),
  @(17,44412,160),
  @(18,44503,160),
  @(19,44491,160),
  @(20,44558,160),
  @(21,44372,160),
  @(22,44448,160),
  @(23,44550,160),
  @(24,44410,120),
  @(25,44561,180),
  @(26,44327,190),
  @(27,44371,180),
  @(28,44585,160),
  @(29,44575,120),
  @(30,44364,160),
  @(31,44239,120),
  @(32,44571,190),
  @(33,44273,160),
  @(34,44414,160),
  @(35,44293,160),
  @(36,44419,130),
  @(37,44484,160),
  @(38,44161,180),
  @(39,44309,160),
  @(40,44266,120),
  @(41,44517,160),
  @(42,44218,130),
  @(43,44447,160),
  @(44,44286,160),
  @(45,44477,160),
  @(46,44523,160),
  @(47,44468,180),
  @(48,44165,68),
  @(49,44512,160),
  @(50,44559,172),
  @(51,44379,160),
  @(52,44432,150),
  @(53,44259,120),
  @(54,44251,80),
  @(55,44496,150),
  @(56,44308,160),
  @(57,44476,160),
  @(58,44427,160),
  @(59,44453,130),
  @(60,44545,180),
  @(61,44216,80),
  @(62,44342,260),
  @(63,44434,140),
  @(64,44580,160),
  @(65,44270,120),
  @(66,44391,160),
  @(67,44189,180),
  @(68,44376,160),
  @(69,44307,130),
  @(70,44417,160),
  @(71,44235,160),
  @(72,44209,160),
  @(73,44250,160),
  @(74,44253,120),
  @(75,44181,90),
  @(76,44298,160),
  @(77,44405,160),
  @(78,44211,120),
  @(79,44377,160),
  @(80,44441,190),
  @(81,44495,160),
  @(82,44505,120),
  @(83,44306,160),
  @(84,44587,120),
  @(85,44384,160),
  @(86,44519,160),
  @(87,44383,180),
  @(88,44469,160),
  @(89,44533,160),
  @(90,44508,160),
  @(91,44295,120),
  @(92,44369,180),
  @(93,44594,130),
  @(94,44526,160),
  @(95,44278,130),
  @(96,44435,810),
  @(97,44314,160),
  @(98,44176,80),
  @(99,44406,160),
  @(100,44467,160),
  @(101,44203,120),
  @(102,44341,160),
  @(103,44245,120),
  @(104,44586,160),
  @(105,44263,180),
  @(106,44305,180),
  @(107,44532,160),
  @(108,44510,160),
  @(109,44442,180),
  @(110,44540,180),
  @(111,44474,160),
  @(112,44350,160),
  @(113,44246,160),
  @(114,44323,160),
  @(115,44398,160),
  @(116,44392,160),
  @(117,44328,160),
  @(118,44589,150),
  @(119,44509,160),
  @(120,44592,160),
  @(121,44518,160),
  @(122,44433,180),
  @(123,44529,160),
  @(124,44382,160),
  @(125,44265,120),
  @(126,44546,180),
  @(127,44363,130),
  @(128,44316,160),
  @(129,44582,180),
  @(130,44159,120),
  @(131,44386,160),
  @(132,44466,160),
  @(133,44322,130),
  @(134,44588,180),
  @(135,44461,160),
  @(136,44452,190),
  @(137,44320,160),
  @(138,44202,120),
  @(139,44560,180),
  @(140,44475,160),
  @(141,44204,180),
  @(142,44358,160),
  @(143,44313,130),
  @(144,44460,160),
  @(145,44333,120),
  @(146,44244,110),
  @(147,44194,80),
  @(148,44172,110),
  @(149,44356,160),
  @(150,44302,130),
  @(151,44539,160),
  @(152,44547,160),
  @(153,44186,180),
  @(154,44482,160),
  @(155,44210,120),
  @(156,44291,89),
  @(157,44217,120),
  @(158,44326,120),
  @(159,44238,130),
  @(160,44188,180),
  @(161,44348,160),
  @(162,44175,120),
  @(163,44201,120),
  @(164,44579,160),
  @(165,44515,160),
  @(166,44438,160),
  @(167,44249,160),
  @(168,44566,130),
  @(169,44351,160),
  @(170,44494,190),
  @(171,44196,180),
  @(172,44237,130),
  @(173,44195,180),
  @(174,44300,160),
  @(175,44578,250),
  @(176,44411,120),
  @(177,44454,160),
  @(178,44281,160),
  @(179,44252,160),
  @(180,44271,180),
  @(181,44420,160),
  @(182,44554,120),
  @(183,44162,160),
  @(184,44343,180),
  @(185,44516,150),
  @(186,44315,130),
  @(187,44568,160),
  @(188,44511,160),
  @(189,44336,160),
  @(190,44231,120),
  @(191,44565,180),
  @(192,44400,160),
  @(193,44334,190),
  @(194,44193,120),
  @(195,44573,160),
  @(196,44403,180),
  @(197,44319,190),
  @(198,44280,120),
  @(199,44362,180),
  @(200,44431,180),
  @(201,44365,180),
  @(202,44567,180),
  @(203,44426,160),
  @(204,44473,160),
  @(205,44357,160),
  @(206,44455,160),
  @(207,44581,130),
  @(208,44553,150),
  @(209,44490,160),
  @(210,44397,160),
  @(211,44446,180),
  @(212,44421,180),
  @(213,44329,160),
  @(214,44208,160),
  @(215,44445,180),
  @(216,44524,160),
  @(217,44355,180),
  @(218,44530,120),
  @(219,44483,180),
  @(220,44294,180),
  @(221,44557,80),
  @(222,44489,160),
  @(223,44264,120),
  @(224,44396,160),
  @(225,44232,120),
  @(226,44279,160),
  @(227,44330,160),
  @(228,44504,160),
  @(229,44572,160),
  @(230,44257,120),
  @(231,44301,130),
  @(232,44370,180),
  @(233,44487,160),
  @(234,44174,180),
  @(235,44200,120),
  @(236,44385,180),
  @(237,44236,120),
  @(238,44413,160),
  @(239,44272,160),
  @(240,44229,160),
  @(241,44214,110),
  @(242,44299,130),
  @(243,44312,160),
  @(244,44399,120),
  @(245,44522,160),
  @(246,44543,160),
  @(247,44167,150),
  @(248,44277,160),
  @(249,44258,230),
  @(250,44349,160),
  @(251,44285,160),
  @(252,44498,160),
  @(253,44179,48),
  @(254,44418,150),
  @(255,44335,160),
  @(256,44552,180),
  @(257,44544,160)
)
foreach ($row in $djData) {
    $ws.Cells.Item($row[0], 4).Value = $row[1]
    $ws.Cells.Item($row[0], 10).Value = $row[2]
}

# Rows whose "Precio minimo/maximo/promedio" (K/L/M) and "Precio $/Kg" (P)
# must change because the special-priced record moved to/from them.
$klmpData = @(
  @(47,1500,1500,1500,500),
  @(48,2000,2000,2000,667),
  @(50,1500,2000,1747,582),
  @(64,1500,1500,1500,500),
  @(82,1500,1500,1500,500),
  @(112,1500,1500,1500,500),
  @(156,1800,1800,1800,600),
  @(177,1500,1500,1500,500),
  @(194,1800,1800,1800,600),
  @(253,2000,2000,2000,667)
)
foreach ($row in $klmpData) {
    $ws.Cells.Item($row[0], 11).Value = $row[1]
    $ws.Cells.Item($row[0], 12).Value = $row[2]
    $ws.Cells.Item($row[0], 13).Value = $row[3]
    $ws.Cells.Item($row[0], 16).Value = $row[4]
}

# Append the brand-new record as row 258 and extend the used range.
$newRow = 258
$ws.Cells.Item($newRow, 1).Value = 3
$ws.Cells.Item($newRow, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($newRow, 3).Value = "Coquimbo"
$ws.Cells.Item($newRow, 4).Value = 44160
$ws.Cells.Item($newRow, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($newRow, 5).Value = 5
$ws.Cells.Item($newRow, 6).Value = 100112039
$ws.Cells.Item($newRow, 7).Value = "Ciboulette"
$ws.Cells.Item($newRow, 8).Value = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value = "Primera"
$ws.Cells.Item($newRow, 10).Value = 230
$ws.Cells.Item($newRow, 11).Value = 1500
$ws.Cells.Item($newRow, 12).Value = 1500
$ws.Cells.Item($newRow, 13).Value = 1500
$ws.Cells.Item($newRow, 14).Value = "`$/docena de atados"
$ws.Cells.Item($newRow, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($newRow, 16).Value = 500
$ws.Cells.Item($newRow, 17).Value = 3
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
